$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C, rows 2 through 344 hold the "Förändrad" (Changed) date.
# Update the serial date value from 45172 to 45175 for every row.
$ws.Range("C2:C344").Value = 45175
